# Apply the cryptos-list price/volume refresh (GitHub Actions data pull).
# Coin names/links are unchanged except rows 44/45, whose coin identities
# swap (VeChain <-> EnergySwap). Numeric-looking "Price" strings are written
# with a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr storage) instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.945.20'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.220.04'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '''291.73'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = '''30.51'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').Value = '''0.0782'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '''49.98'
$ws.Range('E12').Value = '  +5.02%  '
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = '''6.44'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '2.561.62'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').Value = '''13.80'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').Value = '2.224.65'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = '39.878.11'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D21').Value = '''11.10'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').Value = '''5.74'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').Value = '''65.68'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '''237.44'
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('D29').Value = '''9.23'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('E30').Value = '  -6.54%  '
$ws.Range('D31').Value = '''156.65'
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('D32').Value = '''32.03'
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '''2.98'
$ws.Range('E35').Value = '  +6.67%  '
$ws.Range('E36').Value = '  -1.06%  '
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').Value = '''0.0992'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  +1.48%  '
$ws.Range('D41').Value = '''15.34'
$ws.Range('E41').Value = '  -3.41%  '
$ws.Range('D42').Value = '2.108.75'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').Value = '''3.72'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''18.09'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0270'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('D46').Value = '''9.88'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('D47').Value = '''2.01'
$ws.Range('E47').Value = '  -4.34%  '
$ws.Range('D48').Value = '''2.71'
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('D49').Value = '2.433.53'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').Value = '''1.47'
$ws.Range('E50').Value = '  +2.90%  '
$ws.Range('E51').Value = '  +1.83%  '
